$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ---- Sheet: ALC ----
$ws_ALC.Range("H4").Value = 266.66666
$ws_ALC.Range("I4").Value = 266.66666
$ws_ALC.Range("K4").Value = 266.66666
$ws_ALC.Range("M4").Value = -152.66666
$ws_ALC.Range("H49").Value = 1635.375
$ws_ALC.Range("I49").Value = 1436.6
$ws_ALC.Range("J49").Value = 1966.6666
$ws_ALC.Range("K49").Value = 4309.799999999999
$ws_ALC.Range("L49").Value = 5899.9998
$ws_ALC.Range("M49").Value = -4173.799999999999
$ws_ALC.Range("N49").Value = -6171.9998
$ws_ALC.Range("H62").Value = 100911.1
$ws_ALC.Range("I62").Value = 100911.1
$ws_ALC.Range("K62").Value = 100911.1
$ws_ALC.Range("M62").Value = -100287.1
$ws_ALC.Range("H65").Value = 100911.1
$ws_ALC.Range("I65").Value = 100911.1
$ws_ALC.Range("K65").Value = 504555.5
$ws_ALC.Range("M65").Value = -501435.5
$ws_ALC.Range("H80").Value = 785.90625
$ws_ALC.Range("J80").Value = 412.9
$ws_ALC.Range("L80").Value = 1238.7
$ws_ALC.Range("N80").Value = -3234.7
$ws_ALC.Range("H83").Value = 785.90625
$ws_ALC.Range("J83").Value = 412.9
$ws_ALC.Range("L83").Value = 3716.1
$ws_ALC.Range("N83").Value = -13700.1
$ws_ALC.Range("H86").Value = 4360.2144
$ws_ALC.Range("J86").Value = 4934.9
$ws_ALC.Range("L86").Value = 4934.9
$ws_ALC.Range("N86").Value = -7180.9
$ws_ALC.Range("H89").Value = 4360.2144
$ws_ALC.Range("J89").Value = 4934.9
$ws_ALC.Range("L89").Value = 24674.5
$ws_ALC.Range("N89").Value = -35906.5
$ws_ALC.Range("H129").Value = 1296.2424
$ws_ALC.Range("J129").Value = 1714.238
$ws_ALC.Range("L129").Value = 5142.714
$ws_ALC.Range("N129").Value = -15142.714
$ws_ALC.Range("H141").Value = 2549.3044
$ws_ALC.Range("I141").Value = 2078.1428
$ws_ALC.Range("J141").Value = 7496.5
$ws_ALC.Range("K141").Value = 6234.428400000001
$ws_ALC.Range("L141").Value = 22489.5
$ws_ALC.Range("M141").Value = -1054.428400000001
$ws_ALC.Range("N141").Value = -32849.5

# ---- Sheet: ARM ----
$ws_ARM.Range("H45").Value = 2485.5557
$ws_ARM.Range("I45").Value = 2124.3572
$ws_ARM.Range("J45").Value = 3749.75
$ws_ARM.Range("K45").Value = 2124.3572
$ws_ARM.Range("L45").Value = 3749.75
$ws_ARM.Range("M45").Value = -1747.3572
$ws_ARM.Range("N45").Value = -4503.75
$ws_ARM.Range("H61").Value = 11735.692
$ws_ARM.Range("I61").Value = 7509.5713
$ws_ARM.Range("K61").Value = 7509.5713
$ws_ARM.Range("M61").Value = -7297.5713
$ws_ARM.Range("H63").Value = 179313.5
$ws_ARM.Range("I63").Value = 5503
$ws_ARM.Range("J63").Value = 283599.8
$ws_ARM.Range("K63").Value = 5503
$ws_ARM.Range("L63").Value = 283599.8
$ws_ARM.Range("M63").Value = -4817
$ws_ARM.Range("N63").Value = -284971.8
$ws_ARM.Range("H66").Value = 179313.5
$ws_ARM.Range("I66").Value = 5503
$ws_ARM.Range("J66").Value = 283599.8
$ws_ARM.Range("K66").Value = 27515
$ws_ARM.Range("L66").Value = 1417999
$ws_ARM.Range("M66").Value = -24083
$ws_ARM.Range("N66").Value = -1424863
$ws_ARM.Range("H74").Value = 6505.8096
$ws_ARM.Range("I74").Value = 4924.4
$ws_ARM.Range("K74").Value = 4924.4
$ws_ARM.Range("M74").Value = -4050.4
$ws_ARM.Range("H77").Value = 6505.8096
$ws_ARM.Range("I77").Value = 4924.4
$ws_ARM.Range("K77").Value = 24622
$ws_ARM.Range("M77").Value = -20254
$ws_ARM.Range("H132").Value = 5571.8965
$ws_ARM.Range("I132").Value = 4147.4287
$ws_ARM.Range("K132").Value = 12442.2861
$ws_ARM.Range("M132").Value = -9912.286100000001
$ws_ARM.Range("H136").Value = 11735.692
$ws_ARM.Range("I136").Value = 7509.5713
$ws_ARM.Range("K136").Value = 22528.7139
$ws_ARM.Range("M136").Value = -19978.7139

# ---- Sheet: BSM ----
$ws_BSM.Range("H134").Value = 3136.6
$ws_BSM.Range("I134").Value = 3098.5715
$ws_BSM.Range("K134").Value = 9295.7145
$ws_BSM.Range("M134").Value = -6760.7145

# ---- Sheet: CRP ----
$ws_CRP.Range("H3").Value = 9167
$ws_CRP.Range("I3").Value = 12000.5
$ws_CRP.Range("J3").Value = 3500
$ws_CRP.Range("K3").Value = 12000.5
$ws_CRP.Range("L3").Value = 3500
$ws_CRP.Range("M3").Value = -11887.5
$ws_CRP.Range("N3").Value = -3726
$ws_CRP.Range("H31").Value = 1481.5358
$ws_CRP.Range("I31").Value = 1119.32
$ws_CRP.Range("K31").Value = 1119.32
$ws_CRP.Range("M31").Value = -824.3199999999999
$ws_CRP.Range("H34").Value = 1481.5358
$ws_CRP.Range("I34").Value = 1119.32
$ws_CRP.Range("K34").Value = 1119.32
$ws_CRP.Range("M34").Value = -917.3199999999999
$ws_CRP.Range("H132").Value = 1892.7368
$ws_CRP.Range("I132").Value = 1221
$ws_CRP.Range("J132").Value = 3348.1667
$ws_CRP.Range("K132").Value = 3663
$ws_CRP.Range("L132").Value = 10044.5001
$ws_CRP.Range("M132").Value = -1133
$ws_CRP.Range("N132").Value = -15104.5001

# ---- Sheet: CUL ----
$ws_CUL.Range("H94").Value = 4999.5
$ws_CUL.Range("I94").Value = 4999.5
$ws_CUL.Range("J94").Value = 0
$ws_CUL.Range("K94").Value = 14998.5
$ws_CUL.Range("L94").Value = 0
$ws_CUL.Range("M94").Value = -14322.5
$ws_CUL.Range("H131").Value = 1786.675
$ws_CUL.Range("J131").Value = 1974.2903
$ws_CUL.Range("L131").Value = 5922.8709
$ws_CUL.Range("N131").Value = -16002.8709
$ws_CUL.Range("H140").Value = 5262.657
$ws_CUL.Range("I140").Value = 5172.593
$ws_CUL.Range("J140").Value = 5566.625
$ws_CUL.Range("K140").Value = 15517.779
$ws_CUL.Range("L140").Value = 16699.875
$ws_CUL.Range("M140").Value = -10337.779
$ws_CUL.Range("N140").Value = -27059.875
$ws_CUL.Range("N94").ClearContents()

# ---- Sheet: GSM ----
$ws_GSM.Range("H70").Value = 9476.352999999999
$ws_GSM.Range("I70").Value = 6700.125
$ws_GSM.Range("J70").Value = 11944.111
$ws_GSM.Range("K70").Value = 6700.125
$ws_GSM.Range("L70").Value = 11944.111
$ws_GSM.Range("M70").Value = -6430.125
$ws_GSM.Range("N70").Value = -12484.111
$ws_GSM.Range("H73").Value = 9476.352999999999
$ws_GSM.Range("I73").Value = 6700.125
$ws_GSM.Range("J73").Value = 11944.111
$ws_GSM.Range("K73").Value = 6700.125
$ws_GSM.Range("L73").Value = 11944.111
$ws_GSM.Range("M73").Value = -5764.125
$ws_GSM.Range("N73").Value = -13816.111
$ws_GSM.Range("H80").Value = 5262.647
$ws_GSM.Range("J80").Value = 6433.4
$ws_GSM.Range("L80").Value = 6433.4
$ws_GSM.Range("N80").Value = -8429.4
$ws_GSM.Range("H83").Value = 5262.647
$ws_GSM.Range("J83").Value = 6433.4
$ws_GSM.Range("L83").Value = 32167
$ws_GSM.Range("N83").Value = -42151
$ws_GSM.Range("H113").Value = 2999
$ws_GSM.Range("I113").Value = 2999
$ws_GSM.Range("J113").Value = 2999
$ws_GSM.Range("K113").Value = 2999
$ws_GSM.Range("L113").Value = 2999
$ws_GSM.Range("M113").Value = -829
$ws_GSM.Range("N113").Value = -7339
$ws_GSM.Range("H126").Value = 27315.223
$ws_GSM.Range("J126").Value = 54600
$ws_GSM.Range("L126").Value = 163800
$ws_GSM.Range("N126").Value = -168740
$ws_GSM.Range("H132").Value = 10166.75
$ws_GSM.Range("I132").Value = 100012
$ws_GSM.Range("K132").Value = 300036
$ws_GSM.Range("M132").Value = -297506

# ---- Sheet: LTW ----
$ws_LTW.Range("H22").Value = 5598.2144
$ws_LTW.Range("I22").Value = 1150
$ws_LTW.Range("J22").Value = 6339.5835
$ws_LTW.Range("K22").Value = 1150
$ws_LTW.Range("L22").Value = 6339.5835
$ws_LTW.Range("M22").Value = -855
$ws_LTW.Range("N22").Value = -6929.5835
$ws_LTW.Range("H27").Value = 5598.2144
$ws_LTW.Range("I27").Value = 1150
$ws_LTW.Range("J27").Value = 6339.5835
$ws_LTW.Range("K27").Value = 1150
$ws_LTW.Range("L27").Value = 6339.5835
$ws_LTW.Range("M27").Value = -1043
$ws_LTW.Range("N27").Value = -6553.5835
$ws_LTW.Range("H40").Value = 9475
$ws_LTW.Range("I40").Value = 9370
$ws_LTW.Range("K40").Value = 9370
$ws_LTW.Range("M40").Value = -9234
$ws_LTW.Range("H64").Value = 71674
$ws_LTW.Range("J64").Value = 71674
$ws_LTW.Range("L64").Value = 71674
$ws_LTW.Range("N64").Value = -72124
$ws_LTW.Range("H67").Value = 71674
$ws_LTW.Range("J67").Value = 71674
$ws_LTW.Range("L67").Value = 71674
$ws_LTW.Range("N67").Value = -73234
$ws_LTW.Range("H68").Value = 1751.8182
$ws_LTW.Range("I68").Value = 2099.5715
$ws_LTW.Range("J68").Value = 1143.25
$ws_LTW.Range("K68").Value = 2099.5715
$ws_LTW.Range("L68").Value = 1143.25
$ws_LTW.Range("M68").Value = -1350.5715
$ws_LTW.Range("N68").Value = -2641.25
$ws_LTW.Range("H71").Value = 1751.8182
$ws_LTW.Range("I71").Value = 2099.5715
$ws_LTW.Range("J71").Value = 1143.25
$ws_LTW.Range("K71").Value = 10497.8575
$ws_LTW.Range("L71").Value = 5716.25
$ws_LTW.Range("M71").Value = -6753.8575
$ws_LTW.Range("N71").Value = -13204.25
$ws_LTW.Range("H82").Value = 8582.214
$ws_LTW.Range("I82").Value = 17729.334
$ws_LTW.Range("J82").Value = 1721.875
$ws_LTW.Range("K82").Value = 17729.334
$ws_LTW.Range("L82").Value = 1721.875
$ws_LTW.Range("M82").Value = -17368.334
$ws_LTW.Range("N82").Value = -2443.875
$ws_LTW.Range("H85").Value = 8582.214
$ws_LTW.Range("I85").Value = 17729.334
$ws_LTW.Range("J85").Value = 1721.875
$ws_LTW.Range("K85").Value = 17729.334
$ws_LTW.Range("L85").Value = 1721.875
$ws_LTW.Range("M85").Value = -16481.334
$ws_LTW.Range("N85").Value = -4217.875
$ws_LTW.Range("H93").Value = 1717.174
$ws_LTW.Range("I93").Value = 1539.2
$ws_LTW.Range("J93").Value = 2050.875
$ws_LTW.Range("K93").Value = 1539.2
$ws_LTW.Range("L93").Value = 2050.875
$ws_LTW.Range("M93").Value = -291.2
$ws_LTW.Range("N93").Value = -4546.875

# ---- Sheet: WVR ----
$ws_WVR.Range("H81").Value = 4663.4546
$ws_WVR.Range("I81").Value = 4922
$ws_WVR.Range("K81").Value = 9844
$ws_WVR.Range("M81").Value = -8783
$ws_WVR.Range("H84").Value = 4663.4546
$ws_WVR.Range("I84").Value = 4922
$ws_WVR.Range("K84").Value = 49220
$ws_WVR.Range("M84").Value = -43916
$ws_WVR.Range("H109").Value = 87499
$ws_WVR.Range("J109").Value = 87499
$ws_WVR.Range("L109").Value = 87499
$ws_WVR.Range("N109").Value = -90273
$ws_WVR.Range("H122").Value = 3650.889
$ws_WVR.Range("I122").Value = 3650.889
$ws_WVR.Range("J122").Value = 0
$ws_WVR.Range("K122").Value = 10952.667
$ws_WVR.Range("L122").Value = 0
$ws_WVR.Range("M122").Value = -8502.667000000001
$ws_WVR.Range("H132").Value = 1267.4286
$ws_WVR.Range("I132").Value = 1312.1666
$ws_WVR.Range("J132").Value = 999
$ws_WVR.Range("K132").Value = 3936.4998
$ws_WVR.Range("L132").Value = 2997
$ws_WVR.Range("M132").Value = -1406.4998
$ws_WVR.Range("N132").Value = -8057
$ws_WVR.Range("H136").Value = 6421.8853
$ws_WVR.Range("I136").Value = 5533.4053
$ws_WVR.Range("K136").Value = 16600.2159
$ws_WVR.Range("M136").Value = -14050.2159
$ws_WVR.Range("N122").ClearContents()
